# Update constraint & bounds logic
# The "Upper" bound (column C) for 15+ STRIPS, Long Corporate, Equity and
# Liquid Alternatives drops from 1.02 to 1 on every sheet (IBT, Pension,
# Retirement). Also refresh the active-cell selection on each sheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("C3").Value = 1
    $ws.Range("C4").Value = 1
    $ws.Range("C6").Value = 1
    $ws.Range("C7").Value = 1
}

$ibt = $wb.Worksheets.Item("IBT")
$ibt.Activate()
[void]$ibt.Range("C1").Select()

$pension = $wb.Worksheets.Item("Pension")
$pension.Activate()
[void]$pension.Range("C2").Select()

$retirement = $wb.Worksheets.Item("Retirement")
$retirement.Activate()
[void]$retirement.Range("C2").Select()

$ibt.Activate()
